$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 31's format (the special "last row, bottom border" style) onto row 23,
# before we delete rows 24-31 (worker 2's block).
$ws.Range("B31:J31").Copy()
$ws.Range("B23").PasteSpecial(-4122)

# Set final displayed values for the (sorted-ascending) period table, rows 16-23.
$ws.Range("E16").Value2 = "2011"
$ws.Range("F16").Value2 = 40000
$ws.Range("E17").Value2 = "2012"
$ws.Range("F17").Value2 = 40000
$ws.Range("E18").Value2 = "2101"
$ws.Range("F18").Value2 = 40000
$ws.Range("E19").Value2 = "2102"
$ws.Range("F19").Value2 = 40000
$ws.Range("E20").Value2 = "2103"
$ws.Range("F20").Value2 = 40000
$ws.Range("E21").Value2 = "2104"
$ws.Range("F21").Value2 = 40000
$ws.Range("E22").Value2 = "2105"
$ws.Range("F22").Value2 = 40000
$ws.Range("E23").Value2 = "2106"
$ws.Range("F23").Value2 = 28000
$ws.Range("G23").Value2 = 1000000

# Remove worker 2 (HILDA ROSA MIRANDA DE HORTA) rows entirely.
$ws.Rows("24:31").Delete()

# Other scalar updates from the diff.
$ws.Range("E11").Value2 = 308000
$ws.Range("C13").Value2 = 1

# Column D width change.
$ws.Columns("D:D").ColumnWidth = 26.6328125
